# Test Data Driving Values
# - Rename the sheet from "Sheet1" to "ValidLogin"
# - Replace the single "akshara" cell with a small 2x2 login test-data table:
#     A1: UserName   B1: Password
#     A2: admin      B2: pointofsale
# - Leave the active selection on B2 (matches the saved sheetView selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ValidLogin"

# Write the row-2 value first so the shared-string table ends up ordered
# admin, UserName, Password, pointofsale - matching the target workbook.
$ws.Range("A2").Value = "admin"
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "pointofsale"

$ws.Range("B2").Select() | Out-Null
